# Replaces the 25 "NN÷N=" expressions in the single table of the document
# with their updated values, cell by cell (addressed positionally by row/col)
# so that duplicate source/target text values elsewhere in the document cannot
# cause the wrong cell to be updated.

function Replace-CellText {
    param($Table, $Row, $Col, $OldText, $NewText)
    $cellRange = $Table.Cell($Row, $Col).Range
    # Re-wrap the cell range via Document.Range(Start, End): the table cell's own
    # Range object does not reliably bound Find.Execute in this runtime, while a
    # Document.Range(Start, End) correctly restricts the search to that span.
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $ok = $scoped.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 1)
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for Row=$Row Col=$Col ($OldText -> $NewText)"
    }
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

Replace-CellText $t 1 1 "43÷9=" "50÷6="
Replace-CellText $t 1 2 "50÷3=" "87÷2="
Replace-CellText $t 1 3 "25÷3=" "90÷8="
Replace-CellText $t 1 4 "97÷4=" "93÷6="
Replace-CellText $t 1 5 "68÷5=" "90÷7="
Replace-CellText $t 5 1 "80÷6=" "92÷3="
Replace-CellText $t 5 2 "34÷9=" "68÷4="
Replace-CellText $t 5 3 "33÷7=" "73÷8="
Replace-CellText $t 5 4 "24÷2=" "11÷6="
Replace-CellText $t 5 5 "94÷3=" "30÷6="
Replace-CellText $t 9 1 "67÷3=" "54÷6="
Replace-CellText $t 9 2 "32÷2=" "86÷2="
Replace-CellText $t 9 3 "32÷5=" "38÷5="
Replace-CellText $t 9 4 "12÷7=" "18÷5="
Replace-CellText $t 9 5 "32÷5=" "94÷8="
Replace-CellText $t 13 1 "76÷3=" "82÷7="
Replace-CellText $t 13 2 "29÷9=" "35÷3="
Replace-CellText $t 13 3 "20÷3=" "76÷6="
Replace-CellText $t 13 4 "57÷5=" "24÷2="
Replace-CellText $t 13 5 "27÷8=" "11÷5="
Replace-CellText $t 17 1 "24÷3=" "67÷2="
Replace-CellText $t 17 2 "88÷8=" "60÷5="
Replace-CellText $t 17 3 "69÷5=" "35÷5="
Replace-CellText $t 17 4 "24÷2=" "92÷5="
Replace-CellText $t 17 5 "68÷4=" "26÷3="

Write-Host "Done updating division problems."
